$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 31 (shifts existing rows 31-140 down to 32-141)
$ws.Rows.Item(31).Insert()

# Populate the new row 31 with the new weekly price record
$ws.Cells.Item(31, 1).Value = 4
$ws.Cells.Item(31, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(31, 3).Value = "Los Lagos"
$ws.Cells.Item(31, 4).Value = [DateTime]"2021-10-12"
$ws.Cells.Item(31, 5).Value = 10
$ws.Cells.Item(31, 6).Value = 100112017
$ws.Cells.Item(31, 7).Value = "Apio"
$ws.Cells.Item(31, 8).Value = "Americana (o)"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 40
$ws.Cells.Item(31, 11).Value = 11000
$ws.Cells.Item(31, 12).Value = 11000
$ws.Cells.Item(31, 13).Value = 11000
$ws.Cells.Item(31, 14).Value = "`$/docena de matas"
$ws.Cells.Item(31, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(31, 16).Value = 1833
$ws.Cells.Item(31, 17).Value = 6
$ws.Cells.Item(31, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date format as the rest of column D
$ws.Cells.Item(31, 4).NumberFormat = $ws.Cells.Item(32, 4).NumberFormat
